# edit.ps1
# Adds a new "2022-Q1" worksheet (before the "总计" / Total sheet) with
# fund-holding data, and records the new quarter in the "总计" summary sheet.

function Set-TextCell {
    param($cell, $val)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet, positioned right before "总计".
#    Cloning an existing quarter sheet (rather than Worksheets.Add()) carries
#    over the sheetPr / outline settings and lets us inherit the same cell
#    styles (header style, index-column style) used by its sibling sheets.
# ---------------------------------------------------------------------------
$totalSheetBefore = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2021-Q2")
$templateSheet.Copy($totalSheetBefore)
$newSheet = $wb.Worksheets.Item("2021-Q2 (2)")
$newSheet.Name = "2022-Q1"

# NOTE: $totalSheetBefore was captured *before* the sheet copy above, and
# worksheet handles in this host track by tab position, not stable identity
# -- after Copy() shifts everything, that handle now resolves to whatever
# sheet sits at the old position (the just-inserted one). Re-resolve "总计"
# by name now that the sheet count/order has settled.
$totalSheet = $wb.Worksheets.Item("总计")

# The template ("2021-Q2") has 18 data rows (rows 2-19); "2022-Q1" only has 17
# (rows 2-18), so drop the extra templated row before filling in real values.
$newSheet.Rows.Item(19).Delete()

# Header row.
Set-TextCell $newSheet.Range("B1") "基金代码"
Set-TextCell $newSheet.Range("C1") "基金名称"
Set-TextCell $newSheet.Range("D1") "基金规模"
Set-TextCell $newSheet.Range("E1") "股票总仓位"
Set-TextCell $newSheet.Range("F1") "仓位占比"
Set-TextCell $newSheet.Range("G1") "持有市值(亿元)"
Set-TextCell $newSheet.Range("H1") "仓位排名"

# Fund holdings data (code, name, scale, total equity position, position
# weight, held market value, position rank) -- one row per holding.
$fundData = @(
  @("000979", "景顺长城沪港深精选股票", "16.46", "82.61", "9.19", "1.5127", 3),
  @("260112", "景顺长城能源基建混合", "16.49", "60.89", "8.14", "1.3423", 2),
  @("009098", "景顺长城价值领航两年持有期混合", "11.67", "75.58", "9.82", "1.1460", 4),
  @("008850", "景顺长城价值稳进三年定期开放灵活配置混合", "17.06", "69.71", "6.70", "1.1430", 3),
  @("008715", "景顺长城价值驱动一年持有期灵活配置混合型证券投资基金", "16.83", "62.03", "4.79", "0.8062", 5),
  @("013233", "华夏中证500指数智选增强A", "39.48", "92.73", "1.35", "0.5330", 9),
  @("008060", "景顺长城价值边际灵活配置混合", "4.93", "80.78", "9.64", "0.4753", 2),
  @("007994", "华夏中证500指数增强A", "31.45", "92.72", "1.42", "0.4466", 5),
  @("012708", "东方红中证东方红红利低波动指数A", "6.06", "92.33", "1.73", "0.1048", 7),
  @("007995", "华夏中证500指数增强C", "5.45", "92.72", "1.42", "0.0774", 5),
  @("013234", "华夏中证500指数智选增强C", "4.28", "92.73", "1.35", "0.0578", 9),
  @("501219", "华夏智胜先锋股票（LOF）A", "3.61", "94.50", "1.13", "0.0408", 4),
  @("012709", "东方红中证东方红红利低波动指数C", "2.19", "92.33", "1.73", "0.0379", 7),
  @("014198", "华夏智胜先锋股票（LOF）C", "1.30", "94.50", "1.13", "0.0147", 4),
  @("590007", "中邮中证500指数增强A", "0.43", "91.51", "1.38", "0.0059", 10),
  @("008124", "中邮中证500指数增强C", "0.04", "91.51", "1.38", "0.0006", 10),
  @("006992", "嘉合锦创优势精选混合", "0.02", "74.79", "2.42", "0.0005", 7)
)

for ($i = 0; $i -lt $fundData.Count; $i++) {
    $row = $fundData[$i]
    $r = $i + 2
    $newSheet.Range("A$r").Value = $i
    Set-TextCell $newSheet.Range("B$r") $row[0]
    Set-TextCell $newSheet.Range("C$r") $row[1]
    Set-TextCell $newSheet.Range("D$r") $row[2]
    Set-TextCell $newSheet.Range("E$r") $row[3]
    Set-TextCell $newSheet.Range("F$r") $row[4]
    Set-TextCell $newSheet.Range("G$r") $row[5]
    $newSheet.Range("H$r").Value = $row[6]
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" (Total) sheet: insert the new "2022-Q1" row at
#    the top of the data (row 2) and push the existing quarters down.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
# Re-apply the data-row formatting (index-column style etc.) that Insert()
# does not clone, by copying it down from the row just below.
$totalSheet.Range("A3:D3").Copy($totalSheet.Range("A2"))

$totalData = @(
  @("2022-Q1", 17, 7.75),
  @("2021-Q4", 10, 7.66),
  @("2021-Q3", 12, 7.52),
  @("2021-Q2", 18, 7.11),
  @("2021-Q1", 10, 6.19),
  @("2020-Q4", 10, 5.12)
)

for ($i = 0; $i -lt $totalData.Count; $i++) {
    $row = $totalData[$i]
    $r = $i + 2
    $totalSheet.Range("A$r").Value = $i
    Set-TextCell $totalSheet.Range("B$r") $row[0]
    $totalSheet.Range("C$r").Value = $row[1]
    $totalSheet.Range("D$r").Value = $row[2]
}

# Keep the originally active tab selected ("2020-Q4"); adding/copying sheets
# shifts Excel's active-sheet pointer to the most recently touched sheet.
$wb.Worksheets.Item("2020-Q4").Activate()
